$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    $r = $ws.Range($Cell)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '39.959.49'
Set-TextValue 'E2' '  -0.24%  '
Set-TextValue 'D3' '2.204.82'
Set-TextValue 'E3' '  -0.75%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '294.07'
Set-TextValue 'E5' '  +1.36%  '
Set-TextValue 'D6' '87.02'
Set-TextValue 'E6' '  -1.43%  '
Set-TextValue 'D7' '0.512'
Set-TextValue 'E7' '  -0.05%  '
Set-TextValue 'E8' '  +0.03%  '
Set-TextValue 'D9' '0.469'
Set-TextValue 'E9' '  -0.76%  '
Set-TextValue 'B10' 'OKB'
Set-TextValue 'C10' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D10' '51.36'
Set-TextValue 'E10' '  +7.10%  '
Set-TextValue 'B11' 'Avalanche'
Set-TextValue 'C11' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D11' '30.57'
Set-TextValue 'E11' '  +0.08%  '
Set-TextValue 'D12' '0.0779'
Set-TextValue 'E12' '  -0.08%  '
Set-TextValue 'E13' '  +2.65%  '
Set-TextValue 'D14' '6.37'
Set-TextValue 'E14' '  -2.08%  '
Set-TextValue 'D15' '2.551.80'
Set-TextValue 'E15' '  -0.56%  '
Set-TextValue 'D16' '13.77'
Set-TextValue 'E16' '  -1.72%  '
Set-TextValue 'D17' '2.214.20'
Set-TextValue 'E17' '  -0.20%  '
Set-TextValue 'D18' '0.733'
Set-TextValue 'E18' '  +0.28%  '
Set-TextValue 'D19' '39.898.03'
Set-TextValue 'E19' '  -0.26%  '
Set-TextValue 'D20' '0.0₃0886'
Set-TextValue 'E20' '  +0.07%  '
Set-TextValue 'D21' '11.17'
Set-TextValue 'E21' '  -3.88%  '
Set-TextValue 'D22' '5.74'
Set-TextValue 'E22' '  -1.31%  '
Set-TextValue 'D23' '65.36'
Set-TextValue 'E23' '  -0.52%  '
Set-TextValue 'D24' '234.81'
Set-TextValue 'E24' '  -0.32%  '
Set-TextValue 'E25' '  -0.07%  '
Set-TextValue 'D26' '2.46'
Set-TextValue 'E26' '  +0.27%  '
Set-TextValue 'D27' '1.80'
Set-TextValue 'E27' '  -1.56%  '
Set-TextValue 'D28' '23.03'
Set-TextValue 'E28' '  +1.76%  '
Set-TextValue 'E29' '  -4.53%  '
Set-TextValue 'D30' '9.27'
Set-TextValue 'E30' '  +0.43%  '
Set-TextValue 'D31' '159.37'
Set-TextValue 'E31' '  +2.61%  '
Set-TextValue 'D32' '31.67'
Set-TextValue 'E32' '  -0.69%  '
Set-TextValue 'D33' '1.00'
Set-TextValue 'E33' '  +0.05%  '
Set-TextValue 'B34' 'LidoDAOToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D34' '3.04'
Set-TextValue 'E34' '  +6.31%  '
Set-TextValue 'B35' 'Filecoin'
Set-TextValue 'C35' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D35' '4.93'
Set-TextValue 'E35' '  -0.45%  '
Set-TextValue 'D36' '0.0710'
Set-TextValue 'E36' '  -1.23%  '
Set-TextValue 'E37' '  -1.21%  '
Set-TextValue 'E38' '  +1.10%  '
Set-TextValue 'D39' '0.0998'
Set-TextValue 'E39' '  +1.20%  '
Set-TextValue 'D40' '1.74'
Set-TextValue 'E40' '  +2.25%  '
Set-TextValue 'D41' '15.45'
Set-TextValue 'E41' '  -2.41%  '
Set-TextValue 'D42' '2.066.97'
Set-TextValue 'E42' '  -1.96%  '
Set-TextValue 'E43' '  -2.86%  '
Set-TextValue 'D44' '19.33'
Set-TextValue 'E44' '  +9.77%  '
Set-TextValue 'D45' '0.0269'
Set-TextValue 'E45' '  +0.33%  '
Set-TextValue 'D46' '9.86'
Set-TextValue 'E46' '  -0.62%  '
Set-TextValue 'D47' '2.74'
Set-TextValue 'E47' '  +2.56%  '
Set-TextValue 'E48' '  -9.07%  '
Set-TextValue 'D49' '2.423.98'
Set-TextValue 'E49' '  -0.32%  '
Set-TextValue 'E50' '  +1.58%  '
Set-TextValue 'E51' '  +0.22%  '
